# Apply the data-refresh edits to the eICU_24 table (results/table1/eICU_24.docx).
# All edits are addressed by table row/column so they are unambiguous, and the
# full Range.Text of each target cell is rewritten in one shot (rather than via
# Find/Replace on a sub-string) so the existing xml:space="preserve" runs stay
# intact on save.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nbsp = [char]0x00A0
$prefix = "$nbsp$nbsp"

# --- Admission type row (Non-Medical / Medical admission) ---
$t.Cell(18, 2).Range.Text = "24,121 (24.6%)"
$t.Cell(18, 3).Range.Text = "1,336 (17.3%)"
$t.Cell(19, 2).Range.Text = "74,060 (75.4%)"
$t.Cell(19, 3).Range.Text = "6,390 (82.7%)"

# --- Chronic kidney disease labels: capitalize "normal"/"stageN" ---
$t.Cell(39, 1).Range.Text = $prefix + "Normal"
$t.Cell(40, 1).Range.Text = $prefix + "Stage3"
$t.Cell(41, 1).Range.Text = $prefix + "Stage4"
$t.Cell(42, 1).Range.Text = $prefix + "Stage5"

# Row-height bump (616 -> 621 twips = 30.8 -> 31.05 pt) for the stage3/4/5 rows.
$t.Rows.Item(40).Height = 31.05
$t.Rows.Item(41).Height = 31.05
$t.Rows.Item(42).Height = 31.05

# --- SOFA-CNS and MV at 24 hours ---
$t.Cell(44, 2).Range.Text = "0 (0%)"
$t.Cell(44, 3).Range.Text = "0 (0%)"
$t.Cell(44, 4).Range.Text = "NA"
$t.Cell(45, 2).Range.Text = "13,982 (14.2%)"
$t.Cell(45, 3).Range.Text = "3,415 (44.2%)"
$t.Cell(46, 2).Range.Text = "84,199 (85.8%)"
$t.Cell(46, 3).Range.Text = "4,311 (55.8%)"

# --- SOFA - Respiration at 24 hours ---
$t.Cell(48, 2).Range.Text = "746 (0.8%)"
$t.Cell(48, 3).Range.Text = "253 (3.3%)"
$t.Cell(49, 2).Range.Text = "97,435 (99.2%)"
$t.Cell(49, 3).Range.Text = "7,473 (96.7%)"

# --- SOFA - Coagulation at 24 hours ---
$t.Cell(51, 2).Range.Text = "1,930 (2.0%)"
$t.Cell(51, 3).Range.Text = "560 (7.2%)"
$t.Cell(52, 2).Range.Text = "96,251 (98.0%)"
$t.Cell(52, 3).Range.Text = "7,166 (92.8%)"

# --- SOFA - Liver at 24 hours ---
$t.Cell(54, 2).Range.Text = "880 (0.9%)"
$t.Cell(54, 3).Range.Text = "295 (3.8%)"
$t.Cell(55, 2).Range.Text = "97,301 (99.1%)"
$t.Cell(55, 3).Range.Text = "7,431 (96.2%)"

# --- SOFA - Cardiovascular at 24 hours ---
$t.Cell(57, 2).Range.Text = "6,250 (6.4%)"
$t.Cell(57, 3).Range.Text = "1,782 (23.1%)"
$t.Cell(58, 2).Range.Text = "91,931 (93.6%)"
$t.Cell(58, 3).Range.Text = "5,944 (76.9%)"

# --- SOFA - Renal at 24 hours ---
$t.Cell(60, 2).Range.Text = "15,241 (15.5%)"
$t.Cell(60, 3).Range.Text = "2,350 (30.4%)"
$t.Cell(61, 2).Range.Text = "82,940 (84.5%)"
$t.Cell(61, 3).Range.Text = "5,376 (69.6%)"

Write-Host "Done."
